$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes all existing data rows down by one),
# shifting formulas along with it.
$ws.Rows("2:2").Insert(-4121)

# The inserted row picks up formatting from the header row above by default;
# reset it back to the plain "Normal" style used by the rest of the data rows.
$ws.Range("A2:I2").Style = "Normal"

# Column A holds dates formatted like the other rows (re-apply the short date
# format so it reuses the workbook's existing date style).
$ws.Range("A2").NumberFormat = "m/d/yy"
$ws.Range("A2").Value2 = 45950

# New round played at Masterton Golf Course.
$ws.Range("B2").Value2 = "Masterton Golf Course"
$ws.Range("C2").Value2 = "Russell"
$ws.Range("D2").Value2 = "Front-9"
$ws.Range("E2").Value2 = "Practice"
$ws.Range("F2").Formula = "=SUM(6+6+5+8+7+5+6+6)"
$ws.Range("G2").Value2 = 32
$ws.Range("H2").Formula = "=SUM(F2-G2)"
$ws.Range("I2").Value2 = "Only played 8, adjusted score and par for course to account"

# Widen column B to fit the new, longer course name.
$ws.Columns("B:B").ColumnWidth = 20.3

# Match the saved selection/active cell from the source workbook.
$ws.Range("D18").Select()
